$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (attendance count) figures
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 155
$ws1.Range("F4").Value = 256
$ws1.Range("F5").Value = 3961
$ws1.Range("F6").Value = 33
$ws1.Range("F7").Value = 443

# Sheet "全部类型" - same underlying events, different row positions
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 155
$ws4.Range("F4").Value = 256
$ws4.Range("F5").Value = 3961
$ws4.Range("F8").Value = 33
$ws4.Range("F9").Value = 443

$wb.Save()
